# "error rectified for all routes" - append two more rows to the order
# sheet that repeat the existing "Pepper-Rice with Chicken / kilimanjaro /
# coca-cola / emmanuel olajumoke" order, extending the used range from
# A1:D4 to A1:D6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Pepper-Rice with Chicken", "kilimanjaro", "coca-cola", "emmanuel olajumoke"),
    @("Pepper-Rice with Chicken", "kilimanjaro", "coca-cola", "emmanuel olajumoke")
)

$startRow = 5
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
